# "foreach si structuri repetitive" - a new week of attendance (săpt. 11,
# column M) was recorded. For every affected student row the old "săpt. 7"
# (column I) mark is cleared and the later weekly marks are shifted one
# column to the right (I->J->K->L->M), freeing up the new column M for the
# latest week's mark.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row: cells to fully clear (so the XML node disappears, matching the
# saved file) and cells to (re)write with their final value.
$rowPlan = @(
    @{ Row = 6;  Clear = @("I6");  Set = @{ K6 = 2; M6 = 1 } }
    @{ Row = 7;  Clear = @("J7");  Set = @{ K7 = 1; L7 = 2; M7 = 1 } }
    @{ Row = 8;  Clear = @("I8");  Set = @{ J8 = 2; M8 = 1 } }
    @{ Row = 9;  Clear = @("I9");  Set = @{ J9 = 2 } }
    @{ Row = 10; Clear = @("I10"); Set = @{ J10 = 2; K10 = 1; L10 = 2; M10 = 1 } }
    @{ Row = 11; Clear = @("I11"); Set = @{ J11 = 2 } }
    @{ Row = 13; Clear = @("I13"); Set = @{ L13 = 2; M13 = 1 } }
    @{ Row = 14; Clear = @("I14"); Set = @{ K14 = 2; M14 = 1 } }
    @{ Row = 15; Clear = @("K15"); Set = @{ L15 = 2; M15 = 1 } }
    @{ Row = 19; Clear = @("J19"); Set = @{ K19 = 1 } }
    @{ Row = 20; Clear = @("J20"); Set = @{ K20 = 1; L20 = 2; M20 = 1 } }
    @{ Row = 22; Clear = @("I22"); Set = @{ L22 = 2; M22 = 1 } }
)

foreach ($entry in $rowPlan) {
    foreach ($addr in $entry.Clear) {
        $ws.Range($addr).Clear()
    }
    foreach ($addr in $entry.Set.Keys) {
        $ws.Range($addr).Value = $entry.Set[$addr]
    }
}

# The saved workbook view now has the bottom-right pane's selection on M8.
$ws.Range("M8").Select()
